$wb = $excel.ActiveWorkbook

$sheetNames = @("Presentaciones", "Tiempos", "Comidas", "Saludos")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Rename the "Respuesta Incorrecta N" headers to "Opción N"
    $ws.Range("C1").Value = "Opción 1"
    $ws.Range("D1").Value = "Opción 2"
    $ws.Range("E1").Value = "Opción 3"

    # Move the active selection to E2 on every sheet
    $ws.Range("E2").Select() | Out-Null
}

# Re-select the first (active) sheet so it stays the visible/active tab
$wb.Worksheets.Item("Presentaciones").Activate() | Out-Null
$wb.Worksheets.Item("Presentaciones").Range("E2").Select() | Out-Null
